$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a brand-new Sheet2 positioned right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Header row (row 1) -- same header labels as Sheet1, reordered to
# Trip Type / Departure City / Arrival City / Departure Date / Arrival Date / Number of Travellers / Class
$ws2.Range("A1").Value = "Trip Type"
$ws2.Range("B1").Value = "Departure City"
$ws2.Range("C1").Value = "Arrival City"
$ws2.Range("D1").Value = "Departure Date"
$ws2.Range("E1").Value = "Arrival Date"
$ws2.Range("F1").Value = "Number of Travellers"
$ws2.Range("G1").Value = "Class"

# First data row (row 2) -- copy of Sheet1's first data row
$ws2.Range("A2").Value = "One Way"
$ws2.Range("B2").Value = "Kolkata"
$ws2.Range("C2").Value = "Pune"
$ws2.Range("D2").Value = 44110
$ws2.Range("E2").Value = 44110
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = "Economy"

# Formatting: bold header row, wrapped + left/centered text throughout
$ws2.Range("A1:G1").Font.Bold = $true
$ws2.Range("A1:G1").WrapText = $true
$ws2.Range("A1:G1").HorizontalAlignment = -4131
$ws2.Range("A1:G1").VerticalAlignment = -4108

$ws2.Range("A2:G2").WrapText = $true
$ws2.Range("A2:G2").HorizontalAlignment = -4131
$ws2.Range("A2:G2").VerticalAlignment = -4108

# Date columns keep the short-date display used on Sheet1
$ws2.Range("D1:E2").NumberFormat = "m/d/yy"

# Column widths sized to the (smaller) Sheet2 content
$ws2.Columns.Item(1).ColumnWidth = 14.6
$ws2.Columns.Item(2).ColumnWidth = 16.5
$ws2.Columns.Item(3).ColumnWidth = 16.94
$ws2.Columns.Item(4).ColumnWidth = 17.05
$ws2.Columns.Item(5).ColumnWidth = 22.38
$ws2.Columns.Item(6).ColumnWidth = 22.83
$ws2.Columns.Item(7).ColumnWidth = 24.27

# Selection / active-sheet bookkeeping: Sheet1 keeps a range selection and
# loses the active tab, Sheet2 becomes the active/selected sheet.
$ws1.Range("A1:G2").Select()
$ws2.Activate()
$ws2.Range("E6").Select()
